# ----------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
#  - the existing "总计" (grand total) sheet is renamed to "2022-Q1" and its
#    contents are replaced with that quarter's fund-holding detail table
#  - a new "总计" sheet is created (duplicated from the original "总计"
#    sheet, so it keeps the same layout/formatting) containing the original
#    history rows plus a new leading "2022-Q1" row
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-TextCell($rng, $val) {
    # The source numbers in these columns are stored as text (so fixed
    # decimal places such as trailing zeros are not lost). Force the cell to
    # text, assign the value, then drop back to the sheet's default style so
    # no extra number-format styling is left behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$total = $wb.Worksheets.Item("总计")

# 1) Duplicate the current "总计" sheet. The duplicate is inserted right
#    after the original and will become the new "总计" sheet.
$total.Copy($null, $total)
$newTotal = $wb.Worksheets.Item($total.Index + 1)
$newTotal.Name = "总计 (staging)"

# 2) The original sheet becomes "2022-Q1".
$total.Name = "2022-Q1"
$q1 = $total

# 3) The duplicate becomes the new "总计" sheet.
$newTotal.Name = "总计"

# ---------------------------------------------------------------------------
# Rebuild "2022-Q1" (previously "总计") with the fund holdings table.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$q1.Cells.Clear()

# Recreate the sheet's formatting skeleton (header row style + bold/boxed
# index column down to row 16) from an existing sheet that already has the
# identical layout.
$template.Range("A1:H2").Copy()
$q1.Range("A1:H2").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A3:A16").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Cells.Item(2, 1).Value = 0
Set-TextCell $q1.Cells.Item(2, 2) "501203"
Set-TextCell $q1.Cells.Item(2, 3) "易方达创新未来18个月封闭运作混合A"
Set-TextCell $q1.Cells.Item(2, 4) "76.88"
Set-TextCell $q1.Cells.Item(2, 5) "79.51"
Set-TextCell $q1.Cells.Item(2, 6) "6.37"
Set-TextCell $q1.Cells.Item(2, 7) "4.8973"
$q1.Cells.Item(2, 8).Value = 2

$q1.Cells.Item(3, 1).Value = 1
Set-TextCell $q1.Cells.Item(3, 2) "110013"
Set-TextCell $q1.Cells.Item(3, 3) "易方达科翔混合"
Set-TextCell $q1.Cells.Item(3, 4) "49.44"
Set-TextCell $q1.Cells.Item(3, 5) "87.59"
Set-TextCell $q1.Cells.Item(3, 6) "6.92"
Set-TextCell $q1.Cells.Item(3, 7) "3.4212"
$q1.Cells.Item(3, 8).Value = 2

$q1.Cells.Item(4, 1).Value = 2
Set-TextCell $q1.Cells.Item(4, 2) "110001"
Set-TextCell $q1.Cells.Item(4, 3) "易方达平稳增长混合"
Set-TextCell $q1.Cells.Item(4, 4) "33.39"
Set-TextCell $q1.Cells.Item(4, 5) "60.98"
Set-TextCell $q1.Cells.Item(4, 6) "5.53"
Set-TextCell $q1.Cells.Item(4, 7) "1.8465"
$q1.Cells.Item(4, 8).Value = 1

$q1.Cells.Item(5, 1).Value = 3
Set-TextCell $q1.Cells.Item(5, 2) "009341"
Set-TextCell $q1.Cells.Item(5, 3) "易方达均衡成长股票"
Set-TextCell $q1.Cells.Item(5, 4) "64.76"
Set-TextCell $q1.Cells.Item(5, 5) "82.97"
Set-TextCell $q1.Cells.Item(5, 6) "2.71"
Set-TextCell $q1.Cells.Item(5, 7) "1.7550"
$q1.Cells.Item(5, 8).Value = 8

$q1.Cells.Item(6, 1).Value = 4
Set-TextCell $q1.Cells.Item(6, 2) "001018"
Set-TextCell $q1.Cells.Item(6, 3) "易方达新经济灵活配置混合"
Set-TextCell $q1.Cells.Item(6, 4) "65.01"
Set-TextCell $q1.Cells.Item(6, 5) "86.83"
Set-TextCell $q1.Cells.Item(6, 6) "2.48"
Set-TextCell $q1.Cells.Item(6, 7) "1.6122"
$q1.Cells.Item(6, 8).Value = 10

$q1.Cells.Item(7, 1).Value = 5
Set-TextCell $q1.Cells.Item(7, 2) "110029"
Set-TextCell $q1.Cells.Item(7, 3) "易方达科讯混合"
Set-TextCell $q1.Cells.Item(7, 4) "36.09"
Set-TextCell $q1.Cells.Item(7, 5) "91.34"
Set-TextCell $q1.Cells.Item(7, 6) "2.48"
Set-TextCell $q1.Cells.Item(7, 7) "0.8950"
$q1.Cells.Item(7, 8).Value = 10

$q1.Cells.Item(8, 1).Value = 6
Set-TextCell $q1.Cells.Item(8, 2) "519949"
Set-TextCell $q1.Cells.Item(8, 3) "长信利信灵活配置混合A"
Set-TextCell $q1.Cells.Item(8, 4) "2.92"
Set-TextCell $q1.Cells.Item(8, 5) "50.19"
Set-TextCell $q1.Cells.Item(8, 6) "3.87"
Set-TextCell $q1.Cells.Item(8, 7) "0.1130"
$q1.Cells.Item(8, 8).Value = 3

$q1.Cells.Item(9, 1).Value = 7
Set-TextCell $q1.Cells.Item(9, 2) "007294"
Set-TextCell $q1.Cells.Item(9, 3) "长信利信灵活配置混合E"
Set-TextCell $q1.Cells.Item(9, 4) "2.92"
Set-TextCell $q1.Cells.Item(9, 5) "50.19"
Set-TextCell $q1.Cells.Item(9, 6) "3.87"
Set-TextCell $q1.Cells.Item(9, 7) "0.1130"
$q1.Cells.Item(9, 8).Value = 3

$q1.Cells.Item(10, 1).Value = 8
Set-TextCell $q1.Cells.Item(10, 2) "005305"
Set-TextCell $q1.Cells.Item(10, 3) "长信合利混合A"
Set-TextCell $q1.Cells.Item(10, 4) "1.99"
Set-TextCell $q1.Cells.Item(10, 5) "38.27"
Set-TextCell $q1.Cells.Item(10, 6) "3.62"
Set-TextCell $q1.Cells.Item(10, 7) "0.0720"
$q1.Cells.Item(10, 8).Value = 3

$q1.Cells.Item(11, 1).Value = 9
Set-TextCell $q1.Cells.Item(11, 2) "519969"
Set-TextCell $q1.Cells.Item(11, 3) "长信新利灵活配置混合"
Set-TextCell $q1.Cells.Item(11, 4) "0.60"
Set-TextCell $q1.Cells.Item(11, 5) "53.37"
Set-TextCell $q1.Cells.Item(11, 6) "4.44"
Set-TextCell $q1.Cells.Item(11, 7) "0.0266"
$q1.Cells.Item(11, 8).Value = 3

$q1.Cells.Item(12, 1).Value = 10
Set-TextCell $q1.Cells.Item(12, 2) "004608"
Set-TextCell $q1.Cells.Item(12, 3) "长信乐信灵活配置混合A"
Set-TextCell $q1.Cells.Item(12, 4) "0.59"
Set-TextCell $q1.Cells.Item(12, 5) "40.07"
Set-TextCell $q1.Cells.Item(12, 6) "3.44"
Set-TextCell $q1.Cells.Item(12, 7) "0.0203"
$q1.Cells.Item(12, 8).Value = 4

$q1.Cells.Item(13, 1).Value = 11
Set-TextCell $q1.Cells.Item(13, 2) "516910"
Set-TextCell $q1.Cells.Item(13, 3) "富国中证现代物流交易型开放式指数证券投资基金"
Set-TextCell $q1.Cells.Item(13, 4) "0.42"
Set-TextCell $q1.Cells.Item(13, 5) "98.43"
Set-TextCell $q1.Cells.Item(13, 6) "3.88"
Set-TextCell $q1.Cells.Item(13, 7) "0.0163"
$q1.Cells.Item(13, 8).Value = 10

$q1.Cells.Item(14, 1).Value = 12
Set-TextCell $q1.Cells.Item(14, 2) "007293"
Set-TextCell $q1.Cells.Item(14, 3) "长信利信灵活配置混合C"
Set-TextCell $q1.Cells.Item(14, 4) "0.06"
Set-TextCell $q1.Cells.Item(14, 5) "50.19"
Set-TextCell $q1.Cells.Item(14, 6) "3.87"
Set-TextCell $q1.Cells.Item(14, 7) "0.0023"
$q1.Cells.Item(14, 8).Value = 3

$q1.Cells.Item(15, 1).Value = 13
Set-TextCell $q1.Cells.Item(15, 2) "004609"
Set-TextCell $q1.Cells.Item(15, 3) "长信乐信灵活配置混合C"
Set-TextCell $q1.Cells.Item(15, 4) "0.04"
Set-TextCell $q1.Cells.Item(15, 5) "40.07"
Set-TextCell $q1.Cells.Item(15, 6) "3.44"
Set-TextCell $q1.Cells.Item(15, 7) "0.0014"
$q1.Cells.Item(15, 8).Value = 4

$q1.Cells.Item(16, 1).Value = 14
Set-TextCell $q1.Cells.Item(16, 2) "005306"
Set-TextCell $q1.Cells.Item(16, 3) "长信合利混合C"
Set-TextCell $q1.Cells.Item(16, 4) "0.00"
Set-TextCell $q1.Cells.Item(16, 5) "38.27"
Set-TextCell $q1.Cells.Item(16, 6) "3.62"
$q1.Cells.Item(16, 7).Value = 0
$q1.Cells.Item(16, 8).Value = 3

# ---------------------------------------------------------------------------
# Rebuild the new "总计" sheet: same history rows as before, with a new
# "2022-Q1" row inserted at the top (index 0), and the remaining rows'
# running index shifted down by one.
# ---------------------------------------------------------------------------
$rows = @(
    @(0, "2022-Q1", 15, 14.79),
    @(1, "2021-Q4", 23, 15.73),
    @(2, "2021-Q3", 25, 16),
    @(3, "2021-Q2", 18, 17.57),
    @(4, "2021-Q1", 10, 15.16),
    @(5, "2020-Q4", 11, 6.47)
)

# Extend the bold/boxed index-column formatting down to the new last row (7).
$newTotal.Range("A2").Copy()
$newTotal.Range("A7").PasteSpecial(-4122)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $newTotal.Cells.Item($r, 1).Value = $rows[$i][0]
    $newTotal.Cells.Item($r, 2).Value = $rows[$i][1]
    $newTotal.Cells.Item($r, 3).Value = $rows[$i][2]
    $newTotal.Cells.Item($r, 4).Value = $rows[$i][3]
}

# Restore the first sheet as the active/selected tab (duplicating a sheet
# otherwise leaves the new copy selected).
$wb.Worksheets.Item(1).Activate()
